# Add "Area" tracking columns (G, H) and mirrored totals (J, K) to the
# discharge-measurement sheet, per commit "adding area to dishcarge files".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header labels (row 1) ------------------------------------------
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# --- Row 2: summary / total cells ----------------------------------------
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# --- D3:D9 mid-point-of-depth formula, now written as one shared range ---
# (same formula text as before, just re-applied across the whole range so
# it is stored as a single shared formula group)
$ws.Range("D3:D9").Formula = "=(A3/100+(A4/100-A3/100)/2)"

# --- G3: individual area-segment formula ---------------------------------
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# --- G4:G15: shared area-segment formula, extended past the data rows ----
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# --- Restore the view's active selection ---------------------------------
$ws.Range("H6").Select() | Out-Null
